$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-run output: averaged intensities recomputed after adding spiral schemes ---
# Rows 10-16 (scheme index A=8..14) get new values; rows 17-19 (A=15..17) are newly appended.

# Extend column-A header style (bold, bordered, centered) down into the three new rows
# by copying the formatting already used by the existing index column.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

# Row 10: A10 = scheme index 8
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 1.160794165004746
$ws.Cells.Item(10, 4).Value = 0.8636495835432553
$ws.Cells.Item(10, 5).Value = 1.19234943773164
$ws.Cells.Item(10, 6).Value = 0.8687421159298508
$ws.Cells.Item(10, 7).Value = 1.160794165004746
$ws.Cells.Item(10, 8).Value = 0.8636495835432553
$ws.Cells.Item(10, 9).Value = 1.075637167713381
$ws.Cells.Item(10, 10).Value = 0.9958781954106263
$ws.Cells.Item(10, 11).Value = 0.9521629437015247
$ws.Cells.Item(10, 12).Value = 0.849916701545612
$ws.Cells.Item(10, 13).Value = 1.160794165004746
$ws.Cells.Item(10, 14).Value = 1.027999510637448
$ws.Cells.Item(10, 15).Value = 1.021383825552373
$ws.Cells.Item(10, 16).Value = 0.9948912888225796

# Row 11: A11 = scheme index 9
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 0.85322971272447
$ws.Cells.Item(11, 4).Value = 0.0816297148627851
$ws.Cells.Item(11, 5).Value = 1.842163795694577
$ws.Cells.Item(11, 6).Value = 0.7347324975457457
$ws.Cells.Item(11, 7).Value = 0.85322971272447
$ws.Cells.Item(11, 8).Value = 0.0816297148627851
$ws.Cells.Item(11, 9).Value = 1.530015882695179
$ws.Cells.Item(11, 10).Value = 0.9437119560936483
$ws.Cells.Item(11, 11).Value = 1.030452558297665
$ws.Cells.Item(11, 12).Value = 0.343551470266682
$ws.Cells.Item(11, 13).Value = 0.85322971272447
$ws.Cells.Item(11, 14).Value = 0.9618967552786812
$ws.Cells.Item(11, 15).Value = 0.8779389302068945
$ws.Cells.Item(11, 16).Value = 0.919935948522594

# Row 12: A12 = scheme index 10
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 0.8540659325661062
$ws.Cells.Item(12, 4).Value = 0.08179733540637223
$ws.Cells.Item(12, 5).Value = 1.839436282911211
$ws.Cells.Item(12, 6).Value = 0.735308527782781
$ws.Cells.Item(12, 7).Value = 0.8540659325661062
$ws.Cells.Item(12, 8).Value = 0.08179733540637223
$ws.Cells.Item(12, 9).Value = 1.529582797477758
$ws.Cells.Item(12, 10).Value = 0.9433431519359954
$ws.Cells.Item(12, 11).Value = 1.031226202636011
$ws.Cells.Item(12, 12).Value = 0.3435894545819086
$ws.Cells.Item(12, 13).Value = 0.8540659325661062
$ws.Cells.Item(12, 14).Value = 0.9606168091587916
$ws.Cells.Item(12, 15).Value = 0.8776520196666175
$ws.Cells.Item(12, 16).Value = 0.9197937106622678

# Row 13: A13 = scheme index 11
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 0.8526300138102689
$ws.Cells.Item(13, 4).Value = 0.08161538604354578
$ws.Cells.Item(13, 5).Value = 1.843063226281356
$ws.Cells.Item(13, 6).Value = 0.7343769985397989
$ws.Cells.Item(13, 7).Value = 0.8526300138102689
$ws.Cells.Item(13, 8).Value = 0.08161538604354578
$ws.Cells.Item(13, 9).Value = 1.530526519535791
$ws.Cells.Item(13, 10).Value = 0.9440290070386067
$ws.Cells.Item(13, 11).Value = 1.029817421312035
$ws.Cells.Item(13, 12).Value = 0.343111211725477
$ws.Cells.Item(13, 13).Value = 0.8526300138102689
$ws.Cells.Item(13, 14).Value = 0.9623393061624507
$ws.Cells.Item(13, 15).Value = 0.8779214061687423
$ws.Cells.Item(13, 16).Value = 0.9198962230358599

# Row 14: A14 = scheme index 12
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 1.244548000000001
$ws.Cells.Item(14, 4).Value = 0.112552
$ws.Cells.Item(14, 5).Value = 1.095380000000001
$ws.Cells.Item(14, 6).Value = 0.9895520000000008
$ws.Cells.Item(14, 7).Value = 1.244548000000001
$ws.Cells.Item(14, 8).Value = 0.112552
$ws.Cells.Item(14, 9).Value = 1.226851999999999
$ws.Cells.Item(14, 10).Value = 0.7474279999999991
$ws.Cells.Item(14, 11).Value = 1.413472
$ws.Cells.Item(14, 12).Value = 0.4968719999999999
$ws.Cells.Item(14, 13).Value = 1.244548000000001
$ws.Cells.Item(14, 14).Value = 0.6039660000000004
$ws.Cells.Item(14, 15).Value = 0.8605080000000006
$ws.Cells.Item(14, 16).Value = 0.915832

# Row 15: A15 = scheme index 13
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 1.6
$ws.Cells.Item(15, 4).Value = 0.14
$ws.Cells.Item(15, 5).Value = 0.42
$ws.Cells.Item(15, 6).Value = 1.222525000000001
$ws.Cells.Item(15, 7).Value = 1.6
$ws.Cells.Item(15, 8).Value = 0.14
$ws.Cells.Item(15, 9).Value = 0.9515249999999996
$ws.Cells.Item(15, 10).Value = 0.5697249999999997
$ws.Cells.Item(15, 11).Value = 1.764449999999998
$ws.Cells.Item(15, 12).Value = 0.6373875000000008
$ws.Cells.Item(15, 13).Value = 1.6
$ws.Cells.Item(15, 14).Value = 0.28
$ws.Cells.Item(15, 15).Value = 0.8456312500000003
$ws.Cells.Item(15, 16).Value = 0.9132015624999998

# Row 16: A16 = scheme index 14
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 1.350745845350396
$ws.Cells.Item(16, 4).Value = 0.4883842641920018
$ws.Cells.Item(16, 5).Value = 0.6664359366656012
$ws.Cells.Item(16, 6).Value = 1.121165841510399
$ws.Cells.Item(16, 7).Value = 1.350745845350396
$ws.Cells.Item(16, 8).Value = 0.4883842641920018
$ws.Cells.Item(16, 9).Value = 0.9727305456639951
$ws.Cells.Item(16, 10).Value = 0.7550823915520009
$ws.Cells.Item(16, 11).Value = 1.429076401254398
$ws.Cells.Item(16, 12).Value = 0.7836487631872007
$ws.Cells.Item(16, 13).Value = 1.350665731891197
$ws.Cells.Item(16, 14).Value = 0.5774101004288015
$ws.Cells.Item(16, 15).Value = 0.9066829719295992
$ws.Cells.Item(16, 16).Value = 0.945908748671999

# Row 17: A17 = scheme index 15
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.988698126269986
$ws.Cells.Item(17, 4).Value = 0.9934107449639451
$ws.Cells.Item(17, 5).Value = 0.9975992194798802
$ws.Cells.Item(17, 6).Value = 0.9908120491052244
$ws.Cells.Item(17, 7).Value = 0.988698126269986
$ws.Cells.Item(17, 8).Value = 0.9934107449639451
$ws.Cells.Item(17, 9).Value = 0.9939544707871102
$ws.Cells.Item(17, 10).Value = 0.9956131854538658
$ws.Cells.Item(17, 11).Value = 0.9918297630143958
$ws.Cells.Item(17, 12).Value = 0.9908765495072024
$ws.Cells.Item(17, 13).Value = 0.9886723630412961
$ws.Cells.Item(17, 14).Value = 0.9955049822219126
$ws.Cells.Item(17, 15).Value = 0.9926300349547589
$ws.Cells.Item(17, 16).Value = 0.9928492635727012

# Row 18: A18 = scheme index 16
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 1.061426167319558
$ws.Cells.Item(18, 4).Value = 1.045251713196472
$ws.Cells.Item(18, 5).Value = 0.9206179041505751
$ws.Cells.Item(18, 6).Value = 1.027131350457392
$ws.Cells.Item(18, 7).Value = 1.061426167319558
$ws.Cells.Item(18, 8).Value = 1.045251713196472
$ws.Cells.Item(18, 9).Value = 0.949404948704693
$ws.Cells.Item(18, 10).Value = 0.9601305858911097
$ws.Cells.Item(18, 11).Value = 1.015687202764475
$ws.Cells.Item(18, 12).Value = 1.035547049129866
$ws.Cells.Item(18, 13).Value = 1.061426167319558
$ws.Cells.Item(18, 14).Value = 0.9829348086735236
$ws.Cells.Item(18, 15).Value = 1.013606783780999
$ws.Cells.Item(18, 16).Value = 1.001899615201768

# Row 19: A19 = scheme index 17
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 0.9257927643229713
$ws.Cells.Item(19, 4).Value = 1.044620041971311
$ws.Cells.Item(19, 5).Value = 0.9575516246045411
$ws.Cells.Item(19, 6).Value = 1.025591157596936
$ws.Cells.Item(19, 7).Value = 0.9257927643229713
$ws.Cells.Item(19, 8).Value = 1.044620041971311
$ws.Cells.Item(19, 9).Value = 0.9497536781052853
$ws.Cells.Item(19, 10).Value = 1.022324550769879
$ws.Cells.Item(19, 11).Value = 0.9781464310132604
$ws.Cells.Item(19, 12).Value = 1.052356668898478
$ws.Cells.Item(19, 13).Value = 0.9256325400147358
$ws.Cells.Item(19, 14).Value = 1.001085833287926
$ws.Cells.Item(19, 15).Value = 0.9883888971239398
$ws.Cells.Item(19, 16).Value = 0.9945171146603328
